$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New registration row (row 2).
#
# Every cell in the source row is stored as plain text (t="str"), even
# values that look like numbers ("2", "8018389108") or booleans
# ("true"). A leading apostrophe forces Excel to store a value as text
# instead of auto-detecting it as a number/boolean; the style is reset
# back to "Normal" immediately afterwards so the quote-prefix marker
# doesn't leave a lingering cell style behind.
$ws.Range("A2").Value = "kjnkjnnj"
$ws.Range("B2").Value = "sefegf"
$ws.Range("C2").Value = "gfegfe"
$ws.Range("D2").Value = "gfegfeg"
$ws.Range("E2").Value = "aditya@gmail.com"
$ws.Range("F2").Value = "Duo"
$ws.Range("G2").Value = "sfwe"

$ws.Range("H2").Value = "'2"
$ws.Range("H2").Style = "Normal"

$ws.Range("I2").Value = "'true"
$ws.Range("I2").Style = "Normal"

$ws.Range("J2").Value = "GPay"

$ws.Range("K2").Value = "'8018389108"
$ws.Range("K2").Style = "Normal"

$ws.Range("L2").Value = "1769031913866-Gemini_Generated_Image_gmm0m4gmm0m4gmm0.png"

$ws.Range("M2").Value = "'true"
$ws.Range("M2").Style = "Normal"

$ws.Range("N2").Value = "'true"
$ws.Range("N2").Style = "Normal"

$ws.Range("O2").Value = "22/01/2026, 03:15:14"

# The source sheet view is explicitly left-to-right.
$ws.DisplayRightToLeft = $false
